# Generate Report for Handback
# The handback for 6489d606-5d3b-409c-bc93-e42c01cf99eb.md has completed:
# flip its status from "Ready for handoff" to "Handed back: in sync with en-US"
# on every sheet, and record the fresh handback timestamps / clear the stale
# "version not latest" error now that the handback is current.

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusDone
$overview.Range("F3").Value = $statusDone

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusDone
$zhcn.Range("K3").Value = "2016-09-01 06:55:37"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusDone
$dede.Range("K3").Value = "2016-09-01 06:55:44"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
